# Daily update at 8 AM UTC
# Append the next day's tally row (45764 -> 2025-04-17) to the "Wins Over
# Time" sheet, and move the "latest day" short-date formatting from the old
# last row (23) onto the new last row (24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 is no longer the newest entry, so it reverts to the regular
# timestamp format used by every other historical row.
$ws.Range("A23").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 24: today's counts.
$ws.Range("A24").Value = 45764
$ws.Range("B24").Value = 94
$ws.Range("C24").Value = 95
$ws.Range("D24").Value = 95

# Row 24 is now the newest entry, so it gets the short date-only format.
$ws.Range("A24").NumberFormat = "YYYY-MM-DD"
